# ---------------------------------------------------------------------------
# Add 2022-Q4 data:
#  1. Insert a new "2022-Q4" sheet right after "总计", before "2022-Q3".
#  2. Populate it with the fund-holdings table for 2022-Q4.
#  3. Insert a new top data row in "总计" for 2022-Q4 (9 holdings, 0.26亿元)
#     and shift the existing history rows down by one, renumbering the
#     leading index column (A) back to a contiguous 0..7 sequence.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- helper: write a literal (non-numeric-coerced) string into a cell -----
function Set-TextCell($cell, [string]$text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
}

# ============================================================
# 1) Create the new "2022-Q4" worksheet right after "总计"
# ============================================================
$zongji = $wb.Worksheets.Item("总计")
$newWs = $wb.Worksheets.Add($null, $zongji)
$newWs.Name = "2022-Q4"
# NOTE: fetch this reference *after* the insert above - sheet references
# captured before a Worksheets.Add() can end up pointing at stale indices
# once the sheet collection shifts.
$q3 = $wb.Worksheets.Item("2022-Q3")

# Header row (row 1), columns B..H, bold/centered style copied from the
# equivalent header cells on the "2022-Q3" sheet.
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = 2 + $i   # column B=2 .. H=8
    $src = $q3.Cells.Item(1, $col)
    $dst = $newWs.Cells.Item(1, $col)
    $src.Copy($dst)
    $dst.Value = $headers[$i]
}

# Data rows (2022-Q4 fund holdings)
$rows = @(
    @("002317", "招商睿逸稳健配置混合", "4.84", "49.57", "2.69", "0.1302", 8),
    @("217002", "招商安泰平衡混合", "5.50", "49.54", "1.39", "0.0764", 10),
    @("159617", "华夏中证智选500价值稳健策略ETF", "1.74", "97.09", "1.40", "0.0244", 7),
    @("013759", "招商精选平衡混合A", "0.42", "55.24", "3.32", "0.0139", 7),
    @("000646", "华润元大量化优选混合A", "0.19", "65.16", "6.45", "0.0123", 5),
    @("015225", "汇添富中证细分化工产业主题指数增强A", "0.12", "92.36", "2.70", "0.0032", 10),
    @("013760", "招商精选平衡混合C", "0.09", "55.24", "3.32", "0.0030", 7),
    @("015226", "汇添富中证细分化工产业主题指数增强C", "0.03", "92.36", "2.70", "0.0008", 10),
    @("007827", "华润元大量化优选混合C", "0.01", "65.16", "6.45", "0.0006", 5)
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $rowNum = 2 + $r
    $vals = $rows[$r]

    # Column A: numeric index, styled like the other quarter sheets.
    $srcA = $q3.Cells.Item(2, 1)
    $dstA = $newWs.Cells.Item($rowNum, 1)
    $srcA.Copy($dstA)
    $dstA.Value = $r

    Set-TextCell $newWs.Cells.Item($rowNum, 2) $vals[0]     # 基金代码
    Set-TextCell $newWs.Cells.Item($rowNum, 3) $vals[1]     # 基金名称
    Set-TextCell $newWs.Cells.Item($rowNum, 4) $vals[2]     # 基金规模
    Set-TextCell $newWs.Cells.Item($rowNum, 5) $vals[3]     # 股票总仓位
    Set-TextCell $newWs.Cells.Item($rowNum, 6) $vals[4]     # 仓位占比
    Set-TextCell $newWs.Cells.Item($rowNum, 7) $vals[5]     # 持有市值(亿元)
    $newWs.Cells.Item($rowNum, 8).Value = $vals[6]          # 仓位排名 (numeric)
}

# ============================================================
# 2) Update the "总计" summary sheet with the 2022-Q4 row
# ============================================================
# New 2022-Q4 entry followed by the existing history, shifted down by one
# row (the leading index column A is renumbered to stay 0..7 contiguous).
$allRows = @(
    @("2022-Q4", 9, 0.26),
    @("2022-Q3", 8, 0.42),
    @("2022-Q2", 11, 0.8100000000000001),
    @("2022-Q1", 1, 0.02),
    @("2021-Q3", 2, 0.07000000000000001),
    @("2021-Q2", 11, 2.89),
    @("2021-Q1", 25, 5.77),
    @("2020-Q4", 2, 0.03)
)

# Re-fetch "总计" by name too (defensive, see note above).
$zongji = $wb.Worksheets.Item("总计")

for ($r = 0; $r -lt $allRows.Length; $r++) {
    $rowNum = 2 + $r
    $vals = $allRows[$r]

    $srcA = $zongji.Cells.Item(2, 1)
    $dstA = $zongji.Cells.Item($rowNum, 1)
    $srcA.Copy($dstA)
    $dstA.Value = $r

    $zongji.Cells.Item($rowNum, 2).Value = $vals[0]
    $zongji.Cells.Item($rowNum, 3).Value = $vals[1]
    $zongji.Cells.Item($rowNum, 4).Value = $vals[2]
}
